# Keypoints section layout tweak:
#  - Rectangle 18 (the small label shape reading "<div> .content-horizontal-center")
#    inside Group 16 on slide 4 ("Keypoints section") grows slightly, and
#  - its text gains " .keypoints__flex-container" so it reads
#    "<div> .content-horizontal-center .keypoints__flex-container"

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)

$group = $s.Shapes.Item(3)          # "Group 16"
$label = $group.GroupItems.Item(2)  # "Rectangle 18"

# Resize (offset/position is unchanged).
# Shape.Width/Height on a grouped shape are expressed in the group's local
# (child) coordinate space, same units as the raw <a:ext> EMU / 12700.
$label.Width = 226.6224
$label.Height = 10.2877

# Extend the text: keep the existing run's text (it stays a single run)
# and append " .", then append two more runs for the rest of the class name.
$tr = $label.TextFrame.TextRange

$existing = $tr.Characters(1, $tr.Length)
$existing.Text = $existing.Text + " ."

$null = $tr.InsertAfter("keypoints__flex")
$null = $tr.InsertAfter("-container")
